$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update G4:G8 to a uniform calibration factor value
$ws.Range("G4").Value = 45.593649999999997
$ws.Range("G5").Value = 45.593649999999997
$ws.Range("G6").Value = 45.593649999999997
$ws.Range("G7").Value = 45.593649999999997
$ws.Range("G8").Value = 45.593649999999997

# Replace hardcoded H9 value with an AVERAGE formula
$ws.Range("H9").Formula = "=AVERAGE(H3:H8)"

# Update selection to reflect the saved view state
$ws.Range("H10").Select()
